$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.240179
$ws.Range("H2").Value = 3.720537
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 114.155417
$ws.Range("N2").Value = 342.466251
$ws.Range("O2").Value = 0.6835107367845005
$ws.Range("P2").Value = 0.6835107367845005
$ws.Range("Q2").Value = 141.573150899643
$ws.Range("R2").Value = 1274.158358096787
$ws.Range("S2").Value = 0.6835107367845005
$ws.Range("T2").Value = 0.6835107367845005

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.240179
$ws.Range("H3").Value = 3.720537
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 35.924535
$ws.Range("N3").Value = 107.773605
$ws.Range("O3").Value = 0.2150997826628812
$ws.Range("P3").Value = 0.2150997826628812
$ws.Range("Q3").Value = 44.552853891765
$ws.Range("R3").Value = 400.975685025885
$ws.Range("S3").Value = 0.2150997826628812
$ws.Range("T3").Value = 0.2150997826628812

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.240179
$ws.Range("H4").Value = 3.720537
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 16.93339666666667
$ws.Range("N4").Value = 50.80019
$ws.Range("O4").Value = 0.1013894805526183
$ws.Range("P4").Value = 0.1013894805526183
$ws.Range("Q4").Value = 21.00044294467
$ws.Range("R4").Value = 189.00398650203
$ws.Range("S4").Value = 0.1013894805526183
$ws.Range("T4").Value = 0.1013894805526183
